$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-45 with the new meta-regression estimates (var/varname labels + mean/se/ci.lb/ci.ub/pval).
$data = @(
    @('rain_scaled', 'varsel', 0.297, 0.072, 0.157, 0.438, 0),
    @('irr_scaled', 'varsel', -0.055, 0.045, -0.144, 0.034, 0.222),
    @('texture', 'clay', 1.154, 0.206, 0.75, 1.559, 0),
    @('texture', 'clay loam', 1.088, 0.234, 0.628, 1.547, 0),
    @('texture', 'loam', 0.661, 0.154, 0.359, 0.963, 0),
    @('texture', 'sandy clay', 0.693, 0.643, -0.567, 1.954, 0.281),
    @('texture', 'sandy clay loam', 1.874, 0.324, 1.238, 2.509, 0),
    @('texture', 'sandy loam', 0.349, 0.186, -0.015, 0.713, 0.06),
    @('texture', 'silt loam', 0.125, 0.387, -0.634, 0.884, 0.747),
    @('texture', 'silty clay loam', 1.222, 0.246, 0.74, 1.705, 0),
    @('texture', 'unknown', 1.111, 0.302, 0.519, 1.703, 0),
    @('clay_scaled', 'varsel', 0.302, 0.059, 0.186, 0.418, 0),
    @('sand_scaled', 'varsel', -0.185, 0.066, -0.315, -0.055, 0.005),
    @('water_management', 'irrigation', 0.051, 0.545, -1.018, 1.12, 0.926),
    @('water_management', 'irrigation and rain', 0.479, 0.255, -0.021, 0.978, 0.061),
    @('water_management', 'rainfed', 0.96, 0.09, 0.784, 1.136, 0),
    @('sbd_scaled', 'varsel', -0.057, 0.103, -0.259, 0.145, 0.582),
    @('sph_scaled', 'varsel', -0.344, 0.087, -0.515, -0.173, 0),
    @('soc_scaled', 'varsel', -0.05, 0.046, -0.14, 0.04, 0.274),
    @('stn_scaled', 'varsel', 0.483, 0.121, 0.246, 0.719, 0),
    @('crop', 'apple', 0.241, 0.242, -0.234, 0.715, 0.32),
    @('crop', 'beans', 0.853, 0.212, 0.437, 1.269, 0),
    @('crop', 'cabbage', 0.051, 0.485, -0.9, 1.002, 0.916),
    @('crop', 'cocoyam', 0.558, 0.348, -0.125, 1.241, 0.109),
    @('crop', 'lettuce', 0.277, 0.484, -0.671, 1.225, 0.567),
    @('crop', 'maize', 0.778, 0.143, 0.498, 1.058, 0),
    @('crop', 'pea', 1.278, 0.458, 0.381, 2.176, 0.005),
    @('crop', 'peanut', 0.697, 0.395, -0.078, 1.472, 0.078),
    @('crop', 'pepper', 3.337, 0.65, 2.064, 4.61, 0),
    @('crop', 'rice', 1.03, 0.143, 0.749, 1.31, 0),
    @('crop', 'tomato', 3.132, 0.626, 1.904, 4.36, 0),
    @('crop', 'wheat', 0.743, 0.354, 0.05, 1.437, 0.036),
    @('crop_type', 'fruit', 0.241, 0.275, -0.298, 0.78, 0.38),
    @('crop_type', 'grain', 0.94, 0.109, 0.726, 1.155, 0),
    @('crop_type', 'legumes', 0.934, 0.194, 0.553, 1.314, 0),
    @('crop_type', 'others', 0.56, 0.394, -0.213, 1.333, 0.155),
    @('crop_type', 'vegetable', 1.388, 0.302, 0.796, 1.981, 0),
    @('n_fer_scaled', 'varsel', -0.165, 0.073, -0.308, -0.022, 0.024),
    @('p_fer_scaled', 'varsel', 0.014, 0.044, -0.073, 0.101, 0.754),
    @('k_fer_scaled', 'varsel', -0.081, 0.042, -0.163, 0.001, 0.053),
    @('bph_scaled', 'varsel', -0.172, 0.056, -0.282, -0.061, 0.002),
    @('btc_scaled', 'varsel', -0.08, 0.041, -0.16, 0, 0.049),
    @('btn_scaled', 'varsel', -0.048, 0.042, -0.131, 0.035, 0.258),
    @('brate_scaled', 'varsel', -0.021, 0.139, -0.294, 0.253, 0.882)
)

# Rows 2..45 get the refreshed estimates; row index in the sheet = array index + 2.
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# The table now ends at row 45 (three trailing rows - formerly 46, 47, 48 - were dropped).
$ws.Range("A46:A48").EntireRow.Delete()
